$d = $word.ActiveDocument

# Locate the existing list item "La configuration des serveurs ..." - it is the
# last bullet of the "Ce qui marche dans notre projet :" list in the current
# document. We append two new bullet items right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*configuration des serveurs*") {
        $target = $p
    }
}

# Insert a new paragraph right after it; InsertParagraphAfter copies the
# paragraph's formatting (list style/numbering/spacing), so the new bullet
# keeps the same look as the rest of the list.
$target.Range.InsertParagraphAfter()
$first = $target.Next()
$first.Range.InsertAfter("L’utilisation de singletons pour les deux configurations (nombre de connexions, etc) ainsi que la liste des fichiers et des serveurs")

# Insert a second new paragraph after the first new bullet for the second item.
$first.Range.InsertParagraphAfter()
$second = $first.Next()
$second.Range.InsertAfter("L’application s’arrête dès lors qu’on appuie sur echap.")

Write-Output "Inserted two new 'ce qui marche' bullet items."
